# Trade #101 closed at 2026-02-17 09:18:01 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.1
$wsSummary.Range("B4").Value = 0.11
$wsSummary.Range("B6").Value = 101
$wsSummary.Range("B7").Value = 42
$wsSummary.Range("B9").Value = 41.58

# ---------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsStrategy.Range("C4").Value = 100.1
$wsStrategy.Range("D4").Value = 101
$wsStrategy.Range("E4").Value = 0.11
$wsStrategy.Range("F4").Value = 0.1
$wsStrategy.Range("G4").Value = 41.58

# ---------------------------------------------------------------
# Append new trade row (#101) to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 102

    $ws.Cells.Item($row, 1).Value = 101

    # Date-looking text must be forced to stay text, otherwise Excel's
    # auto-detection turns "2026-02-17" into a date serial number.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).Value = "09:17:55"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.11
    $ws.Cells.Item($row, 7).Value = 0.128205
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 16.5501
    $ws.Cells.Item($row, 10).Value = 0.02
    $ws.Cells.Item($row, 11).Value = 100.1
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14
}
